# Add 2022-Q4 data.
#
# The workbook tracks quarterly fund-holding snapshots: a "总计" (totals) sheet
# lists one summary row per quarter, and each quarter also gets its own detail
# sheet with per-fund figures. This change adds a new "2022-Q4" quarter:
#   - a new "2022-Q4" detail sheet (a fresh snapshot of the same fund, with
#     updated figures), inserted right after the "总计" sheet and before the
#     existing "2022-Q3" sheet
#   - a new row on "总计" for "2022-Q4", with the older rows following after it

$wb = $excel.ActiveWorkbook

# --- 1. Create the "2022-Q4" detail sheet ------------------------------
# Duplicate "2022-Q3" (same columns/formatting) positioned immediately
# before it, then rename the copy and refresh its figures for the new
# quarter. The fund itself (code/name) is unchanged.
$prevQuarter = $wb.Worksheets.Item("2022-Q3")
$prevQuarter.Copy($prevQuarter)
$newQuarter = $wb.ActiveSheet
$newQuarter.Name = "2022-Q4"

$newQuarter.Range("D2").Value = "'3.11"
$newQuarter.Range("E2").Value = "'92.49"
$newQuarter.Range("F2").Value = "'2.78"
$newQuarter.Range("G2").Value = "'0.0865"
$newQuarter.Range("H2").Value = 4

# --- 2. Update the "总计" summary sheet --------------------------------
# Existing rows 2 & 3 ("2022-Q3" and "2021-Q3") each slide down to make
# room for the new "2022-Q4" row at the top of the data, and a fresh row
# 4 is appended for what used to be the last quarter, "2021-Q3".
$tot = $wb.Worksheets.Item("总计")

$tot.Range("B2").Value = "2022-Q4"

$tot.Range("B3").Value = "2022-Q3"
$tot.Range("D3").Value = 0.09

$tot.Range("A3").Copy($tot.Range("A4"))
$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 1
$tot.Range("D4").Value = 0.01
